$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure every cell in the target range is stored as literal text (matches
# the source workbook's inlineStr / string cells) rather than being
# auto-coerced into numbers by COM when the value looks numeric.
$ws.Range("A1:V7").NumberFormat = "@"

$data = @(
  @('accountName', 'cost', 'impressions', 'clicks', 'averageCpc', 'invalidClicks', 'invalidClickRate', 'impressionShare', 'exactMatchImpressionShare', 'budgetLostImpressionShare', 'qualityLostImpressionShare', 'conversions', 'convRate', 'convValue', 'costPerConv', 'valuePerConv', 'allConv', 'allConvRate', 'allConvValue', 'costPerAllConv', 'valuePerAllConv', 'accountid'),
  @('C3sFHc9Txk', '1,584,598', '6,604,331', '13,962,080', '21,156.94', '141,605,592', '21,087.74', '20,590.70', '20,598.08', '20,701.05', '21,098.50', '20,950.55', '20,960.78', '20,764.06', '21,153.01', '20,750.28', '65,418,560.96', '20,880.43', '66,483,063.26', '20,741.32', '21,285.33', '70,795'),
  @('SoTVSVuH0w', '1,584,598', '6,604,331', '13,962,080', '21,156.94', '141,605,592', '21,087.74', '20,590.70', '20,598.08', '20,701.05', '21,098.50', '20,950.55', '20,960.78', '20,764.06', '21,153.01', '20,750.28', '65,418,560.96', '20,880.43', '66,483,063.26', '20,741.32', '21,285.33', '70,795'),
  @('aqZaakvSNq', '1,584,598', '6,604,331', '13,962,080', '21,156.94', '141,605,592', '21,087.74', '20,590.70', '20,598.08', '20,701.05', '21,098.50', '20,950.55', '20,960.78', '20,764.06', '21,153.01', '20,750.28', '65,418,560.96', '20,880.43', '66,483,063.26', '20,741.32', '21,285.33', '70,795'),
  @('imbbyolD4T', '1,584,598', '6,604,331', '13,962,080', '21,156.94', '141,605,592', '21,087.74', '20,590.70', '20,598.08', '20,701.05', '21,098.50', '20,950.55', '20,960.78', '20,764.06', '21,153.01', '20,750.28', '65,418,560.96', '20,880.43', '66,483,063.26', '20,741.32', '21,285.33', '70,795'),
  @('OEXiADPjQd', '804,765', '3,365,646', '7,464,471', '20,516.57', '73,586,762', '21,136.75', '20,643.88', '21,130.84', '21,024.48', '20,864.22', '21,263.28', '20,890.87', '20,539.14', '20,816.25', '20,560.25', '34,263,346.94', '20,951.21', '33,334,928.96', '21,126.07', '20,480.14', '18,734'),
  @('0WSSd6geVU', '804,765', '3,365,646', '7,464,471', '20,516.57', '73,586,762', '21,136.75', '20,643.88', '21,130.84', '21,024.48', '20,864.22', '21,263.28', '20,890.87', '20,539.14', '20,816.25', '20,560.25', '34,263,346.94', '20,951.21', '33,334,928.96', '21,126.07', '20,480.14', '18,734'),
)

for ($r = 0; $r -lt $data.Length; $r++) {
  $rowValues = $data[$r]
  for ($c = 0; $c -lt $rowValues.Length; $c++) {
    $ws.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
  }
}
